$d = $word.ActiveDocument

# --- 1. Remove the "Meta description: ..." paragraph that follows the H1 title ---
$metaPara = $d.Paragraphs(2)
$null = $metaPara.Range.Delete()

# --- 2. Insert a new bold "Play Druid's Dream ..." paragraph right before the
#        final (image-prompt) paragraph, then turn the final paragraph's text
#        into the "Read our review ..." meta-description text, keeping the
#        italic formatting that paragraph already has. ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)

# Create a fresh empty paragraph right before the last one.
$null = $lastPara.Range.InsertParagraphBefore()

# That new paragraph is now the second-to-last paragraph; fill it with clean
# OOXML so the run layout matches the rest of the document (leading empty
# run + bold text run).
$newCount = $d.Paragraphs.Count
$newPara = $d.Paragraphs($newCount - 1)
$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Druid''s Dream For Free - Stunning Nature-Themed Slot Machine</w:t></w:r></w:p>'
$null = $newPara.Range.InsertXML($newXml)

# Swap the old image-generation-prompt text in the final paragraph for the
# meta-description copy, preserving the paragraph's existing (italic) run
# formatting. Trim the trailing paragraph mark off the Range before setting
# .Text so the replacement happens in place (instead of being inserted as a
# new, differently formatted run) and so autocorrect "smart quotes" aren't
# applied to the straight apostrophes in the replacement text.
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$textRange = $finalPara.Range.Duplicate()
$null = $textRange.MoveEnd(1, -1)
$textRange.Text = "Read our review of Druid's Dream and play for free. A captivating nature-themed slot machine with stunning graphics and a wide betting range."
